$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 130
$ws.Range("B2").Value = 2048
$ws.Range("C2").Value = 1024
$ws.Range("D2").Value = 8
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 1024
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 3712360
$ws.Range("I2").Value = 0.275835
$ws.Range("J2").Value = 564.91
$ws.Range("K2").Value = 29002.7
$ws.Range("L2").Value = 29032.1
$ws.Range("M2").Value = 28934.4
$ws.Range("N2").Value = 29025.2
$ws.Range("O2").Value = 29020.5
$ws.Range("P2").Value = 29009.7
$ws.Range("Q2").Value = 77.401
$ws.Range("R2").Value = 95.2076
$ws.Range("S2").Value = 47.7616
$ws.Range("T2").Value = 94.797
$ws.Range("U2").Value = 93.9948
$ws.Range("V2").Value = 63.1103
$ws.Range("W2").Value = 14.1306
$ws.Range("X2").Value = 14.153
$ws.Range("Y2").Value = 14.0892
$ws.Range("Z2").Value = 14.1495
$ws.Range("AA2").Value = 14.146
$ws.Range("AB2").Value = 14.1304
$ws.Range("AC2").Value = 0
$ws.Range("AD2").Value = 100
$ws.Range("AE2").Value = 97.20871845005333
$ws.Range("AF2").Value = 0.0006955505630481
$ws.Range("AG2").Value = 95.8941650263266
$ws.Range("AH2").Value = 94.62826195070886
